# Update countries & provincias Spain
# Refresh COVID-19 case counters to the 25-Sep-2020 14:11 snapshot, which
# re-ranks a handful of countries by total cases (column B, descending),
# and bump the "datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Estados Unidos (row 4)
$ws.Range("B4").Value = 7187179
$ws.Range("C4").Value = 1708
$ws.Range("D4").Value = 4438906
$ws.Range("E4").Value = 2540718
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 207555

# Kuwait (row 41)
$ws.Range("B41").Value = 102441
$ws.Range("C41").Value = 590
$ws.Range("D41").Value = 93562
$ws.Range("E41").Value = 8284
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 595

# Nepal overtakes Venezuela and Costa Rica in the ranking (rows 53-55)
$ws.Range("A53").Value = "Nepal"
$ws.Range("B53").Value = 70614
$ws.Range("C53").Value = 1313
$ws.Range("D53").Value = 51866
$ws.Range("E53").Value = 18289
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 459

$ws.Range("A54").Value = "Venezuela"
$ws.Range("B54").Value = 70406
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 59745
$ws.Range("E54").Value = 10080
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 581

$ws.Range("A55").Value = "Costa Rica"
$ws.Range("B55").Value = 69459
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 26554
$ws.Range("E55").Value = 42110
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 795

# Azerbaiyan (row 68)
$ws.Range("B68").Value = 39787
$ws.Range("C68").Value = 101
$ws.Range("D68").Value = 37392
$ws.Range("E68").Value = 1812
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 583

# Estado de Palestina (row 70)
$ws.Range("B70").Value = 37963
$ws.Range("C70").Value = 372
$ws.Range("D70").Value = 27183
$ws.Range("E70").Value = 10502
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 278

# Australia (row 78)
$ws.Range("D78").Value = 24523
$ws.Range("E78").Value = 1608

# Dinamarca (row 80)
$ws.Range("B80").Value = 25594
$ws.Range("C80").Value = 678
$ws.Range("D80").Value = 19010
$ws.Range("E80").Value = 5937
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 647

# Madagascar (row 88)
$ws.Range("B88").Value = 16221
$ws.Range("C88").Value = 30
$ws.Range("D88").Value = 14867
$ws.Range("E88").Value = 1126
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 228

# Finlandia (row 104)
$ws.Range("B104").Value = 9484
$ws.Range("C104").Value = 105
$ws.Range("E104").Value = 1291

# Uganda (row 113)
$ws.Range("B113").Value = 7218
$ws.Range("C113").Value = 154
$ws.Range("D113").Value = 3611
$ws.Range("E113").Value = 3536
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 71

# Islandia (row 149)
$ws.Range("B149").Value = 2561
$ws.Range("C149").Value = 49
$ws.Range("D149").Value = 2151
$ws.Range("E149").Value = 400

# Lesoto (row 162)
$ws.Range("B162").Value = 1558
$ws.Range("C162").Value = 4
$ws.Range("D162").Value = 797
$ws.Range("E162").Value = 726

# Vietnam (row 168)
$ws.Range("D168").Value = 999
$ws.Range("E168").Value = 35

# Islas Feroe (row 179)
$ws.Range("B179").Value = 458
$ws.Range("C179").Value = 3
$ws.Range("D179").Value = 417
$ws.Range("E179").Value = 41

# Gibraltar (row 183)
$ws.Range("B183").Value = 364
$ws.Range("C183").Value = 3
$ws.Range("D183").Value = 333
$ws.Range("E183").Value = 31

# Islas Malvinas overtakes Montserrat in the ranking (rows 215-216)
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1

# Bump the "datos actualizados" timestamp shown in row 1
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 14:11"
